$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.507.55'
$ws.Range('E2').Value = '  +4.01%  '
$ws.Range('D3').Value = '3.341.47'
$ws.Range('E3').Value = '  +4.27%  '
$ws.Range('E4').Value = '  -0.03%  '
$r = $ws.Range('D5')
$r.NumberFormat = "@"
$r.Value = '561.75'
$r.Style = "Normal"
$ws.Range('E5').Value = '  +4.60%  '
$r = $ws.Range('D6')
$r.NumberFormat = "@"
$r.Value = '151.91'
$r.Style = "Normal"
$ws.Range('E6').Value = '  +4.64%  '
$r = $ws.Range('D7')
$r.NumberFormat = "@"
$r.Value = '0.999'
$r.Style = "Normal"
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '3.340.53'
$ws.Range('E8').Value = '  +4.02%  '
$r = $ws.Range('D9')
$r.NumberFormat = "@"
$r.Value = '0.536'
$r.Style = "Normal"
$ws.Range('E9').Value = '  +1.28%  '
$r = $ws.Range('D10')
$r.NumberFormat = "@"
$r.Value = '7.41'
$r.Style = "Normal"
$ws.Range('E10').Value = '  +1.10%  '
$ws.Range('E11').Value = '  +3.82%  '
$r = $ws.Range('D12')
$r.NumberFormat = "@"
$r.Value = '0.433'
$r.Style = "Normal"
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').Value = '3.917.03'
$ws.Range('E13').Value = '  +4.32%  '
$ws.Range('E14').Value = '  +0.49%  '
$r = $ws.Range('D15')
$r.NumberFormat = "@"
$r.Value = '26.78'
$r.Style = "Normal"
$ws.Range('E15').Value = '  +3.40%  '
$ws.Range('E16').Value = '  +3.04%  '
$ws.Range('D17').Value = '62.504.20'
$ws.Range('E17').Value = '  +3.95%  '
$ws.Range('D18').Value = '3.337.80'
$ws.Range('E18').Value = '  +4.06%  '
$r = $ws.Range('D19')
$r.NumberFormat = "@"
$r.Value = '6.34'
$r.Style = "Normal"
$ws.Range('E19').Value = '  +1.24%  '
$r = $ws.Range('D20')
$r.NumberFormat = "@"
$r.Value = '13.79'
$r.Style = "Normal"
$ws.Range('E20').Value = '  +4.72%  '
$r = $ws.Range('D21')
$r.NumberFormat = "@"
$r.Value = '8.39'
$r.Style = "Normal"
$ws.Range('E21').Value = '  +1.07%  '
$r = $ws.Range('D22')
$r.NumberFormat = "@"
$r.Value = '385.08'
$r.Style = "Normal"
$ws.Range('E22').Value = '  +2.47%  '
$ws.Range('E23').Value = '  -0.19%  '
$ws.Range('E24').Value = '  +1.53%  '
$r = $ws.Range('D25')
$r.NumberFormat = "@"
$r.Value = '70.03'
$r.Style = "Normal"
$ws.Range('E25').Value = '  +0.00%  '
$r = $ws.Range('D26')
$r.NumberFormat = "@"
$r.Value = '0.179'
$r.Style = "Normal"
$ws.Range('E26').Value = '  +5.55%  '
$r = $ws.Range('D27')
$r.NumberFormat = "@"
$r.Value = '9.04'
$r.Style = "Normal"
$ws.Range('E27').Value = '  +3.16%  '
$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').Value = '0.0₃0950'
$ws.Range('E28').Value = '  +6.05%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$r = $ws.Range('D29')
$r.NumberFormat = "@"
$r.Value = '0.999'
$r.Style = "Normal"
$ws.Range('E29').Value = '  -0.02%  '
$r = $ws.Range('D30')
$r.NumberFormat = "@"
$r.Value = '6.61'
$r.Style = "Normal"
$ws.Range('E30').Value = '  +7.64%  '
$ws.Range('E31').Value = '  +3.96%  '
$r = $ws.Range('D32')
$r.NumberFormat = "@"
$r.Value = '5.58'
$r.Style = "Normal"
$ws.Range('E32').Value = '  +3.88%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$r = $ws.Range('D33')
$r.NumberFormat = "@"
$r.Value = '1.31'
$r.Style = "Normal"
$ws.Range('E33').Value = '  +9.14%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$r = $ws.Range('D34')
$r.NumberFormat = "@"
$r.Value = '22.85'
$r.Style = "Normal"
$ws.Range('E34').Value = '  +2.23%  '
$r = $ws.Range('D35')
$r.NumberFormat = "@"
$r.Value = '6.73'
$r.Style = "Normal"
$ws.Range('E35').Value = '  +0.39%  '
$ws.Range('E36').Value = '  +9.10%  '
$r = $ws.Range('D37')
$r.NumberFormat = "@"
$r.Value = '159.14'
$r.Style = "Normal"
$ws.Range('E37').Value = '  +1.64%  '
$ws.Range('E38').Value = '  +11.85%  '
$r = $ws.Range('D39')
$r.NumberFormat = "@"
$r.Value = '26.88'
$r.Style = "Normal"
$ws.Range('E39').Value = '  +4.69%  '
$ws.Range('E40').Value = '  +4.47%  '
$ws.Range('D41').Value = '2.785.59'
$ws.Range('E41').Value = '  -0.63%  '
$ws.Range('E42').Value = '  +7.07%  '
$r = $ws.Range('D44')
$r.NumberFormat = "@"
$r.Value = '40.46'
$r.Style = "Normal"
$ws.Range('E44').Value = '  +1.46%  '
$ws.Range('E45').Value = '  +3.43%  '
$ws.Range('E46').Value = '  +5.04%  '
$r = $ws.Range('D47')
$r.NumberFormat = "@"
$r.Value = '22.04'
$r.Style = "Normal"
$ws.Range('E47').Value = '  +6.90%  '
$ws.Range('D48').Value = '3.382.25'
$ws.Range('E48').Value = '  +4.19%  '
$ws.Range('E49').Value = '  -2.31%  '
$r = $ws.Range('D50')
$r.NumberFormat = "@"
$r.Value = '6.30'
$r.Style = "Normal"
$ws.Range('E50').Value = '  +2.17%  '
$r = $ws.Range('D51')
$r.NumberFormat = "@"
$r.Value = '287.39'
$r.Style = "Normal"
$ws.Range('E51').Value = '  +7.79%  '
